$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.341.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.441.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.73%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.995.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.388.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.988"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "415.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("E22").Value = "  +5.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "86.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "608.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.43%  "
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.144"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.95%  "
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.223.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -7.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.57%  "

# Row 16/17 swap: TRON <-> WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.442.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.09%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.121"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

# Row 25/26 swap: RenderToken <-> InternetComputer(DFINITY)
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.66%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.93%  "
